# Sync attendance_reports: normalize "Recorded By" (column G) value ordering.
# For every data row, if the G-column value is a comma-separated list whose
# first entry is literally "System", move that "System" token to the end of
# the list (preserving the order of the remaining entries) while leaving any
# value that does not start with "System, " untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }

    if ($val.StartsWith("System, ")) {
        $parts = $val.Split(",")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }
        $rest = $parts[1..($parts.Length - 1)]
        $newParts = $rest + @("System")
        $newVal = [string]::Join(", ", $newParts)
        $cell.Value2 = $newVal
    }
}
